$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BENCHMARK")

# Update J13: Azami 9.335 TL -> Azami 995,5 TL
$ws.Range("J13").Value = "Hesaba: Asgari 1 TL | Azami 995,5 TL"

# Clear F24 and F25 (remove their contents entirely)
$ws.Range("F24").ClearContents()
$ws.Range("F25").ClearContents()
